# Updates crypto price/volume figures (and restores Kaspa/VeChain row order)
# to match the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.069.43"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "3.195.05"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.13"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.195.90"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.66"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.503"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.22%  "
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.24"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "3.719.67"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "66.215.78"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "3.197.44"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "507.72"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.730"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.99"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.58"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.35"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  +39.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.92"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.97"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  -5.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.48"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.43"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "499.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").Value = "0.0₃0768"
$ws.Range("E39").Value = "  +13.88%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0419"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.02"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.71"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.296"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").Value = "2.909.95"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.40"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.23%  "
